$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 769; this shifts rows 769:816 down to 770:817
$ws.Rows("769:769").Insert()

# Populate the newly inserted row 769 with a copy of the (now shifted) row 770's
# data, except for the date in column D which gets a new value.
$ws.Range("A769").Value = 9
$ws.Range("B769").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C769").Value = "Metropolitana"
$ws.Range("D769").Value = 45106
$ws.Range("E769").Value = 13
$ws.Range("F769").Value = 100112024
$ws.Range("G769").Value = "Choclo"
$ws.Range("H769").Value = "Dulce o Americano"
$ws.Range("I769").Value = "Primera"
$ws.Range("J769").Value = 70
$ws.Range("K769").Value = 16000
$ws.Range("L769").Value = 18000
$ws.Range("M769").Value = 17000
$ws.Range("N769").Value = "`$/malla 70 unidades"
$ws.Range("O769").Value = "Región de Arica y Parinacota"
$ws.Range("P769").Value = 243
$ws.Range("Q769").Value = 70
$ws.Range("R769").Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D
$ws.Range("D769").NumberFormat = $ws.Range("D770").NumberFormat
